$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (G=4564)
$ws.Range("H6").Value = 235.16667
$ws.Range("I6").Value = 235.16667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 705.50001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -593.50001
$ws.Range("N6").ClearContents()
# Row 86 (G=12603)
$ws.Range("H86").Value = 2549.5
$ws.Range("I86").Value = 2332.5
$ws.Range("J86").Value = 2875
$ws.Range("K86").Value = 2332.5
$ws.Range("L86").Value = 2875
$ws.Range("M86").Value = -1209.5
$ws.Range("N86").Value = -5121
# Row 89 (G=12603)
$ws.Range("H89").Value = 2549.5
$ws.Range("I89").Value = 2332.5
$ws.Range("J89").Value = 2875
$ws.Range("K89").Value = 11662.5
$ws.Range("L89").Value = 14375
$ws.Range("M89").Value = -6046.5
$ws.Range("N89").Value = -25607
# Row 127 (G=36114)
$ws.Range("H127").Value = 899.75
$ws.Range("I127").Value = 500
$ws.Range("K127").Value = 1500
$ws.Range("M127").Value = 3460
# Row 131 (G=36108)
$ws.Range("H131").Value = 6343.3687
$ws.Range("I131").Value = 1147.3572
$ws.Range("K131").Value = 3442.0716
$ws.Range("M131").Value = 1597.9284

$ws = $wb.Worksheets.Item("ARM")
# Row 11 (G=3767)
$ws.Range("H11").Value = 2525000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
# Row 33 (G=3352)
$ws.Range("H33").Value = 13475
$ws.Range("I33").Value = 13475
$ws.Range("K33").Value = 13475
$ws.Range("M33").Value = -13146
# Row 61 (G=43999)
$ws.Range("H61").Value = 4592
$ws.Range("I61").Value = 5878
$ws.Range("J61").Value = 3857.1428
$ws.Range("K61").Value = 5878
$ws.Range("L61").Value = 3857.1428
$ws.Range("M61").Value = -5666
$ws.Range("N61").Value = -4281.1428
# Row 74 (G=44000)
$ws.Range("H74").Value = 2050.4
$ws.Range("I74").Value = 2050.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2050.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1176.4
$ws.Range("N74").ClearContents()
# Row 77 (G=44000)
$ws.Range("H77").Value = 2050.4
$ws.Range("I77").Value = 2050.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10252
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5884
$ws.Range("N77").ClearContents()
# Row 97 (G=19941)
$ws.Range("H97").Value = 772.875
$ws.Range("I97").Value = 725.5
$ws.Range("J97").Value = 915
$ws.Range("K97").Value = 725.5
$ws.Range("L97").Value = 915
$ws.Range("M97").Value = -229.5
$ws.Range("N97").Value = -1907
# Row 136 (G=43999)
$ws.Range("H136").Value = 4592
$ws.Range("I136").Value = 5878
$ws.Range("J136").Value = 3857.1428
$ws.Range("K136").Value = 17634
$ws.Range("L136").Value = 11571.4284
$ws.Range("M136").Value = -15084
$ws.Range("N136").Value = -16671.4284

$ws = $wb.Worksheets.Item("BSM")
# Row 40 (G=19514)
$ws.Range("H40").Value = 44834
$ws.Range("J40").Value = 44834
$ws.Range("L40").Value = 44834
$ws.Range("N40").Value = -45364
# Row 86 (G=12526)
$ws.Range("H86").Value = 74395.14
$ws.Range("I86").Value = 3318
$ws.Range("J86").Value = 145472.28
$ws.Range("K86").Value = 3318
$ws.Range("L86").Value = 145472.28
$ws.Range("M86").Value = -2195
$ws.Range("N86").Value = -147718.28
# Row 89 (G=12526)
$ws.Range("H89").Value = 74395.14
$ws.Range("I89").Value = 3318
$ws.Range("J89").Value = 145472.28
$ws.Range("K89").Value = 16590
$ws.Range("L89").Value = 727361.4
$ws.Range("M89").Value = -10974
$ws.Range("N89").Value = -738593.4
# Row 96 (G=19525)
$ws.Range("H96").Value = 101978
$ws.Range("I96").Value = 109375.8
$ws.Range("J96").Value = 28000
$ws.Range("K96").Value = 109375.8
$ws.Range("L96").Value = 28000
$ws.Range("M96").Value = -106629.8
$ws.Range("N96").Value = -33492
# Row 127 (G=35358)
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws = $wb.Worksheets.Item("CUL")
# Row 17 (G=4640)
$ws.Range("H17").Value = 5333.3335
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5333.3335
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 16000.0005
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -16338.0005
# Row 39 (G=4712)
$ws.Range("H39").Value = 5457
$ws.Range("J39").Value = 5457
$ws.Range("L39").Value = 16371
$ws.Range("N39").Value = -16959
# Row 58 (G=4703)
$ws.Range("H58").Value = 2353.3333
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 2485.7144
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 7457.1432
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = -7713.1432
# Row 93 (G=19808)
$ws.Range("H93").Value = 1275
$ws.Range("J93").Value = 1350
$ws.Range("L93").Value = 4050
$ws.Range("N93").Value = -7794
# Row 133 (G=44073)
$ws.Range("H133").Value = 4503.7144
$ws.Range("I133").Value = 1756
$ws.Range("J133").Value = 5362.375
$ws.Range("K133").Value = 5268
$ws.Range("L133").Value = 16087.125
$ws.Range("M133").Value = -208
$ws.Range("N133").Value = -26207.125
# Row 134 (G=44074)
$ws.Range("H134").Value = 3595.2964
$ws.Range("I134").Value = 2090.7693
$ws.Range("K134").Value = 6272.3079
$ws.Range("M134").Value = -1202.3079
# Row 137 (G=44088)
$ws.Range("H137").Value = 3379
$ws.Range("J137").Value = 3558.0833
$ws.Range("L137").Value = 10674.2499
$ws.Range("N137").Value = -20874.2499
# Row 139 (G=44102)
$ws.Range("H139").Value = 2054.2258
$ws.Range("I139").Value = 1335.75
$ws.Range("J139").Value = 3360.5454
$ws.Range("K139").Value = 4007.25
$ws.Range("L139").Value = 10081.6362
$ws.Range("M139").Value = 1132.75
$ws.Range("N139").Value = -20361.6362

$ws = $wb.Worksheets.Item("GSM")
# Row 13 (G=2443)
$ws.Range("H13").Value = 16628.111
$ws.Range("I13").Value = 410.6
$ws.Range("J13").Value = 36900
$ws.Range("K13").Value = 410.6
$ws.Range("L13").Value = 36900
$ws.Range("M13").Value = -271.6
$ws.Range("N13").Value = -37178
# Row 18 (G=4309)
$ws.Range("H18").Value = 28500
$ws.Range("J18").Value = 21333.334
$ws.Range("L18").Value = 21333.334
$ws.Range("N18").Value = -21919.334
# Row 80 (G=12521)
$ws.Range("H80").Value = 3223.4211
$ws.Range("I80").Value = 3229.1667
$ws.Range("J80").Value = 3213.5715
$ws.Range("K80").Value = 3229.1667
$ws.Range("L80").Value = 3213.5715
$ws.Range("M80").Value = -2231.1667
$ws.Range("N80").Value = -5209.5715
# Row 83 (G=12521)
$ws.Range("H83").Value = 3223.4211
$ws.Range("I83").Value = 3229.1667
$ws.Range("J83").Value = 3213.5715
$ws.Range("K83").Value = 16145.8335
$ws.Range("L83").Value = 16067.8575
$ws.Range("M83").Value = -11153.8335
$ws.Range("N83").Value = -26051.8575
# Row 109 (G=25691)
$ws.Range("H109").Value = 15642.143
$ws.Range("J109").Value = 15642.143
$ws.Range("L109").Value = 15642.143
$ws.Range("N109").Value = -17722.143
# Row 123 (G=34150)
$ws.Range("H123").Value = 9158.888999999999
$ws.Range("J123").Value = 9158.888999999999
$ws.Range("L123").Value = 9158.888999999999
$ws.Range("N123").Value = -14058.889
# Row 126 (G=36184)
$ws.Range("H126").Value = 3901.2
$ws.Range("I126").Value = 3235.3333
$ws.Range("J126").Value = 4900
$ws.Range("K126").Value = 9705.999899999999
$ws.Range("L126").Value = 14700
$ws.Range("M126").Value = -7235.999899999999
$ws.Range("N126").Value = -19640

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (G=2631)
$ws.Range("H2").Value = 9378
# Row 14 (G=3771)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
# Row 34 (G=3347)
$ws.Range("H34").Value = 16779.8
$ws.Range("I34").Value = 18474.75
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 18474.75
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -18302.75
$ws.Range("N34").Value = -10344
# Row 40 (G=36248)
$ws.Range("H40").Value = 2750
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1364
$ws.Range("N40").Value = -4272
# Row 64 (G=10810)
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450
# Row 67 (G=10810)
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560
# Row 68 (G=12563)
$ws.Range("H68").Value = 2633.5
$ws.Range("I68").Value = 1733.25
$ws.Range("J68").Value = 3833.8333
$ws.Range("K68").Value = 1733.25
$ws.Range("L68").Value = 3833.8333
$ws.Range("M68").Value = -984.25
$ws.Range("N68").Value = -5331.8333
# Row 71 (G=12563)
$ws.Range("H71").Value = 2633.5
$ws.Range("I71").Value = 1733.25
$ws.Range("J71").Value = 3833.8333
$ws.Range("K71").Value = 8666.25
$ws.Range("L71").Value = 19169.1665
$ws.Range("M71").Value = -4922.25
$ws.Range("N71").Value = -26657.1665
# Row 82 (G=12565)
$ws.Range("H82").Value = 2563.7856
$ws.Range("I82").Value = 1765.5555
$ws.Range("K82").Value = 1765.5555
$ws.Range("M82").Value = -1404.5555
# Row 85 (G=12565)
$ws.Range("H85").Value = 2563.7856
$ws.Range("I85").Value = 1765.5555
$ws.Range("K85").Value = 1765.5555
$ws.Range("M85").Value = -517.5554999999999
# Row 100 (G=19995)
$ws.Range("H100").Value = 3725
$ws.Range("I100").Value = 3300
$ws.Range("J100").Value = 3866.6667
$ws.Range("K100").Value = 3300
$ws.Range("L100").Value = 3866.6667
$ws.Range("M100").Value = -2759
$ws.Range("N100").Value = -4948.6667
# Row 109 (G=27209)
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws = $wb.Worksheets.Item("WVR")
# Row 11 (G=3001)
$ws.Range("H11").Value = 17000.6
$ws.Range("I11").Value = 21667.666
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 21667.666
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -21525.666
$ws.Range("N11").Value = -10284
# Row 64 (G=11036)
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# Row 67 (G=11036)
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# Row 96 (G=19977)
$ws.Range("H96").Value = 1519
$ws.Range("I96").Value = 1285.7142
$ws.Range("J96").Value = 1723.125
$ws.Range("K96").Value = 1285.7142
$ws.Range("L96").Value = 1723.125
$ws.Range("M96").Value = 87.28580000000011
$ws.Range("N96").Value = -4469.125
# Row 123 (G=34127)
$ws.Range("H123").Value = 24182.125
$ws.Range("J123").Value = 24182.125
$ws.Range("L123").Value = 24182.125
$ws.Range("N123").Value = -33982.125
# Row 125 (G=34276)
$ws.Range("H125").Value = 58742.6
$ws.Range("J125").Value = 58742.6
$ws.Range("L125").Value = 58742.6
$ws.Range("N125").Value = -68582.60000000001
# Row 135 (G=42043)
$ws.Range("H135").Value = 61016.43
$ws.Range("J135").Value = 61016.43
$ws.Range("L135").Value = 61016.43
$ws.Range("N135").Value = -71156.42999999999
